# Add an "image" column to the product table (column E) and populate
# the example row with a placeholder image URL.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = "image"
$ws.Range("E2").Value = "https://placehold.co/600x400/EEE/31343C"

# Matches the selection recorded in the saved workbook (E4).
$ws.Range("E4").Select()
